# Add an extended command to connect to WiFi on the "MCU-ESP commands" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MCU-ESP commands")

# Insert a new row right after the existing "Get Wifi status of ESP8266" entry
# (row 5), pushing all subsequent rows down by one.
$ws.Rows.Item(6).Insert() | Out-Null

$ws.Cells.Item(6, 1).Value = "84 F0 80 85"
$ws.Cells.Item(6, 2).Value = "MCU -> ESP"
$ws.Cells.Item(6, 3).Value = "Connect to WiFi"

# Reflect the new active selection at the bottom of the table.
$ws.Activate()
$ws.Range("C12").Select() | Out-Null
